$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1967213114754098
$ws.Cells.Item(2, 3).Value = 0.5491803278688525
$ws.Cells.Item(2, 10).Value = 0.01229508196721311
$ws.Cells.Item(2, 16).Value = 0.1680327868852459
$ws.Cells.Item(2, 19).Value = 0.07377049180327869

$ws.Cells.Item(3, 2).Value = 0.0072992700729927
$ws.Cells.Item(3, 3).Value = 0.05109489051094891
$ws.Cells.Item(3, 10).Value = 0.0364963503649635
$ws.Cells.Item(3, 16).Value = 0.7372262773722628
$ws.Cells.Item(3, 19).Value = 0.1678832116788321

$ws.Cells.Item(4, 10).Value = 0.1052631578947368
$ws.Cells.Item(4, 16).Value = 0.6578947368421053
$ws.Cells.Item(4, 19).Value = 0.2368421052631579

$ws.Cells.Item(6, 2).Value = 0.06
$ws.Cells.Item(6, 4).Value = 0.01
$ws.Cells.Item(6, 6).Value = 0.095
$ws.Cells.Item(6, 10).Value = 0.25
$ws.Cells.Item(6, 15).Value = 0.01
$ws.Cells.Item(6, 17).Value = 0.12
$ws.Cells.Item(6, 18).Value = 0.095
$ws.Cells.Item(6, 19).Value = 0.36

$ws.Cells.Item(7, 2).Value = 0.1151079136690648
$ws.Cells.Item(7, 4).Value = 0.01438848920863309
$ws.Cells.Item(7, 6).Value = 0.05755395683453238
$ws.Cells.Item(7, 10).Value = 0.09352517985611511
$ws.Cells.Item(7, 15).Value = 0.01438848920863309
$ws.Cells.Item(7, 17).Value = 0.1294964028776978
$ws.Cells.Item(7, 18).Value = 0.07194244604316546
$ws.Cells.Item(7, 19).Value = 0.5035971223021583

$ws.Cells.Item(8, 2).Value = 0.09022556390977443
$ws.Cells.Item(8, 4).Value = 0.02506265664160401
$ws.Cells.Item(8, 6).Value = 0.06516290726817042
$ws.Cells.Item(8, 10).Value = 0.1553884711779449
$ws.Cells.Item(8, 15).Value = 0.01503759398496241
$ws.Cells.Item(8, 17).Value = 0.1528822055137845
$ws.Cells.Item(8, 18).Value = 0.05764411027568922
$ws.Cells.Item(8, 19).Value = 0.4385964912280702

$ws.Cells.Item(9, 2).Value = 0.07692307692307693
$ws.Cells.Item(9, 4).Value = 0.02797202797202797
$ws.Cells.Item(9, 6).Value = 0.08391608391608392
$ws.Cells.Item(9, 10).Value = 0.1678321678321678
$ws.Cells.Item(9, 15).Value = 0.01398601398601399
$ws.Cells.Item(9, 17).Value = 0.1468531468531468
$ws.Cells.Item(9, 18).Value = 0.0979020979020979
$ws.Cells.Item(9, 19).Value = 0.3846153846153846

$ws.Cells.Item(10, 2).Value = 0.1000841042893188
$ws.Cells.Item(10, 4).Value = 0.0159798149705635
$ws.Cells.Item(10, 5).Value = 0.004205214465937763
$ws.Cells.Item(10, 6).Value = 0.05971404541631623
$ws.Cells.Item(10, 10).Value = 0.144659377628259
$ws.Cells.Item(10, 15).Value = 0.0176619007569386
$ws.Cells.Item(10, 17).Value = 0.1976450798990748
$ws.Cells.Item(10, 18).Value = 0.08662741799831791
$ws.Cells.Item(10, 19).Value = 0.3734230445752734

$ws.Cells.Item(11, 7).Value = 0.141025641025641
$ws.Cells.Item(11, 10).Value = 0.1324786324786325
$ws.Cells.Item(11, 11).Value = 0.1965811965811966
$ws.Cells.Item(11, 12).Value = 0.5085470085470085
$ws.Cells.Item(11, 19).Value = 0.02136752136752137

$ws.Cells.Item(12, 7).Value = 0.7603305785123967
$ws.Cells.Item(12, 10).Value = 0.1735537190082645
$ws.Cells.Item(12, 11).Value = 0.008264462809917356
$ws.Cells.Item(12, 12).Value = 0.01652892561983471
$ws.Cells.Item(12, 19).Value = 0.04132231404958678

$ws.Cells.Item(13, 7).Value = 0.55
$ws.Cells.Item(13, 10).Value = 0.4
$ws.Cells.Item(13, 19).Value = 0.05

$ws.Cells.Item(15, 6).Value = 0.01142857142857143
$ws.Cells.Item(15, 8).Value = 0.2057142857142857
$ws.Cells.Item(15, 9).Value = 0.05142857142857143
$ws.Cells.Item(15, 10).Value = 0.4
$ws.Cells.Item(15, 11).Value = 0.08
$ws.Cells.Item(15, 13).Value = 0.01714285714285714
$ws.Cells.Item(15, 15).Value = 0.04
$ws.Cells.Item(15, 19).Value = 0.1942857142857143

$ws.Cells.Item(16, 6).Value = 0.02958579881656805
$ws.Cells.Item(16, 8).Value = 0.2130177514792899
$ws.Cells.Item(16, 9).Value = 0.1005917159763314
$ws.Cells.Item(16, 10).Value = 0.408284023668639
$ws.Cells.Item(16, 11).Value = 0.0650887573964497
$ws.Cells.Item(16, 13).Value = 0.005917159763313609
$ws.Cells.Item(16, 15).Value = 0.04142011834319527
$ws.Cells.Item(16, 19).Value = 0.136094674556213

$ws.Cells.Item(17, 6).Value = 0.0339943342776204
$ws.Cells.Item(17, 8).Value = 0.1784702549575071
$ws.Cells.Item(17, 9).Value = 0.08781869688385269
$ws.Cells.Item(17, 10).Value = 0.4702549575070821
$ws.Cells.Item(17, 11).Value = 0.0594900849858357
$ws.Cells.Item(17, 13).Value = 0.0226628895184136
$ws.Cells.Item(17, 15).Value = 0.0594900849858357
$ws.Cells.Item(17, 19).Value = 0.08781869688385269

$ws.Cells.Item(18, 6).Value = 0.03571428571428571
$ws.Cells.Item(18, 8).Value = 0.1488095238095238
$ws.Cells.Item(18, 9).Value = 0.09523809523809523
$ws.Cells.Item(18, 10).Value = 0.4345238095238095
$ws.Cells.Item(18, 11).Value = 0.08333333333333333
$ws.Cells.Item(18, 13).Value = 0.01785714285714286
$ws.Cells.Item(18, 15).Value = 0.07738095238095238
$ws.Cells.Item(18, 19).Value = 0.1071428571428571

$ws.Cells.Item(19, 6).Value = 0.01188299817184644
$ws.Cells.Item(19, 8).Value = 0.2138939670932358
$ws.Cells.Item(19, 9).Value = 0.06398537477148081
$ws.Cells.Item(19, 10).Value = 0.3985374771480805
$ws.Cells.Item(19, 11).Value = 0.1087751371115174
$ws.Cells.Item(19, 13).Value = 0.02285191956124314
$ws.Cells.Item(19, 15).Value = 0.06764168190127971
$ws.Cells.Item(19, 19).Value = 0.1124314442413163
